$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59: Raw TA opened crm 10/10/19
$ws.Range("A58").Copy()
$ws.Range("A59").PasteSpecial(-4122)
$ws.Range("A59").Value = 43748
$ws.Range("B59").Value = 2365.542
$ws.Range("C59").Value = 2207.0300000000002
$ws.Range("D59").Formula = "=100*(B59-C59)/C59"
$ws.Range("E59").Value = 169
$ws.Range("F59").Value = "Raw TA opened crm 10/10/19"

# Row 60: TA evap; opened crm 10/10/19
$ws.Range("A58").Copy()
$ws.Range("A60").PasteSpecial(-4122)
$ws.Range("A60").Value = 43748
$ws.Range("B60").Value = 2357.459
$ws.Range("C60").Value = 2207.0300000000002
$ws.Range("D60").Formula = "=100*(B60-C60)/C60"
$ws.Range("E60").Value = 169
$ws.Range("F60").Value = "TA evap; opened crm 10/10/19"

# Update view: top-left cell and selection shift to match new data added
$excel.ActiveWindow.ScrollRow = 39
$ws.Range("F50").Select()
